# Add a new "jumpPower" field/column to the Character sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Character")

# Insert a new column at D. This shifts the existing "weight" column (D) to E,
# and -- because D was previously part of the merged width range C:D (14.5) --
# the insert naturally extends that merged range to C:E, matching the target
# column widths/merge exactly.
$ws.Columns("D:D").Insert()

# The old "weight" column landed in E after the insert; move its contents back
# to D (its style/number format/border, s="1"/s="2", already carried over
# correctly as part of the column insert/shift).
$ws.Range("D1").Value = "weight"
$ws.Range("D4").Value = "float"
$ws.Range("D5").Value = 50

# Populate the real new column (E) with the new jumpPower field.
$ws.Range("E1").Value = "jumpPower"
$ws.Range("E4").Value = "float"
$ws.Range("E5").Value = 100

# moveSpeed data value changed from 10 to 20.
$ws.Range("C5").Value = 20

# Update the active selection to match.
$ws.Range("C5").Select()
